$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 72-73, shifting existing rows 72-138 down to 74-140
$ws.Rows("72:73").Insert()

$ws.Cells.Item(72, 1).Value = 2
$ws.Cells.Item(72, 2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(72, 3).Value = 'Coquimbo'
$ws.Cells.Item(72, 4).Value = 44574
$ws.Cells.Item(72, 5).Value = 4
$ws.Cells.Item(72, 6).Value = 100112031
$ws.Cells.Item(72, 7).Value = 'Poroto verde'
$ws.Cells.Item(72, 8).Value = 'Magnum'
$ws.Cells.Item(72, 9).Value = 'Primera'
$ws.Cells.Item(72, 10).Value = 700
$ws.Cells.Item(72, 11).Value = 18000
$ws.Cells.Item(72, 12).Value = 20000
$ws.Cells.Item(72, 13).Value = 19000
$ws.Cells.Item(72, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(72, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(72, 16).Value = 760
$ws.Cells.Item(72, 17).Value = 25
$ws.Cells.Item(72, 18).Value = 'Hortaliza'
$ws.Cells.Item(73, 1).Value = 2
$ws.Cells.Item(73, 2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(73, 3).Value = 'Coquimbo'
$ws.Cells.Item(73, 4).Value = 44574
$ws.Cells.Item(73, 5).Value = 4
$ws.Cells.Item(73, 6).Value = 100112031
$ws.Cells.Item(73, 7).Value = 'Poroto verde'
$ws.Cells.Item(73, 8).Value = 'Sin especificar'
$ws.Cells.Item(73, 9).Value = 'Primera'
$ws.Cells.Item(73, 10).Value = 400
$ws.Cells.Item(73, 11).Value = 25000
$ws.Cells.Item(73, 12).Value = 27000
$ws.Cells.Item(73, 13).Value = 26000
$ws.Cells.Item(73, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(73, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(73, 16).Value = 1040
$ws.Cells.Item(73, 17).Value = 25
$ws.Cells.Item(73, 18).Value = 'Hortaliza'
